$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 319 (shifts old rows 319+ down to 322+)
$ws.Rows("319:321").Insert()

# Row 319 (new, copy of old row 316 content)
$ws.Cells.Item(319,1).Value = 10
$ws.Cells.Item(319,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(319,3).Value = "La Araucanía"
$ws.Cells.Item(319,4).Value = 44468
$ws.Cells.Item(319,5).Value = 9
$ws.Cells.Item(319,6).Value = "Fruta"
$ws.Cells.Item(319,7).Value = 100102
$ws.Cells.Item(319,8).Value = "Cítricos"
$ws.Cells.Item(319,9).Value = 100102004
$ws.Cells.Item(319,10).Value = "Mandarina"
$ws.Cells.Item(319,11).Value = "Murcott"
$ws.Cells.Item(319,12).Value = "Primera"
$ws.Cells.Item(319,13).Value = 200
$ws.Cells.Item(319,14).Value = 12000
$ws.Cells.Item(319,15).Value = 12000
$ws.Cells.Item(319,16).Value = 12000
$ws.Cells.Item(319,17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(319,18).Value = "Región de O'Higgins"
$ws.Cells.Item(319,19).Value = 667
$ws.Cells.Item(319,20).Value = 18

# Row 320 (new, copy of old row 317 content)
$ws.Cells.Item(320,1).Value = 10
$ws.Cells.Item(320,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(320,3).Value = "La Araucanía"
$ws.Cells.Item(320,4).Value = 44468
$ws.Cells.Item(320,5).Value = 9
$ws.Cells.Item(320,6).Value = "Fruta"
$ws.Cells.Item(320,7).Value = 100102
$ws.Cells.Item(320,8).Value = "Cítricos"
$ws.Cells.Item(320,9).Value = 100102004
$ws.Cells.Item(320,10).Value = "Mandarina"
$ws.Cells.Item(320,11).Value = "Murcott"
$ws.Cells.Item(320,12).Value = "Primera"
$ws.Cells.Item(320,13).Value = 200
$ws.Cells.Item(320,14).Value = 10000
$ws.Cells.Item(320,15).Value = 10000
$ws.Cells.Item(320,16).Value = 10000
$ws.Cells.Item(320,17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(320,18).Value = "Región de O'Higgins"
$ws.Cells.Item(320,19).Value = 500
$ws.Cells.Item(320,20).Value = 20

# Row 321 (new, copy of old row 318 content)
$ws.Cells.Item(321,1).Value = 10
$ws.Cells.Item(321,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(321,3).Value = "La Araucanía"
$ws.Cells.Item(321,4).Value = 44468
$ws.Cells.Item(321,5).Value = 9
$ws.Cells.Item(321,6).Value = "Fruta"
$ws.Cells.Item(321,7).Value = 100102
$ws.Cells.Item(321,8).Value = "Cítricos"
$ws.Cells.Item(321,9).Value = 100102004
$ws.Cells.Item(321,10).Value = "Mandarina"
$ws.Cells.Item(321,11).Value = "Murcott"
$ws.Cells.Item(321,12).Value = "Segunda"
$ws.Cells.Item(321,13).Value = 120
$ws.Cells.Item(321,14).Value = 8000
$ws.Cells.Item(321,15).Value = 8000
$ws.Cells.Item(321,16).Value = 8000
$ws.Cells.Item(321,17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(321,18).Value = "Región de O'Higgins"
$ws.Cells.Item(321,19).Value = 444
$ws.Cells.Item(321,20).Value = 18

# Row 316 (updated: date changed)
$ws.Cells.Item(316,1).Value = 10
$ws.Cells.Item(316,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(316,3).Value = "La Araucanía"
$ws.Cells.Item(316,4).Value = 44476
$ws.Cells.Item(316,5).Value = 9
$ws.Cells.Item(316,6).Value = "Fruta"
$ws.Cells.Item(316,7).Value = 100102
$ws.Cells.Item(316,8).Value = "Cítricos"
$ws.Cells.Item(316,9).Value = 100102004
$ws.Cells.Item(316,10).Value = "Mandarina"
$ws.Cells.Item(316,11).Value = "Murcott"
$ws.Cells.Item(316,12).Value = "Primera"
$ws.Cells.Item(316,13).Value = 200
$ws.Cells.Item(316,14).Value = 12000
$ws.Cells.Item(316,15).Value = 12000
$ws.Cells.Item(316,16).Value = 12000
$ws.Cells.Item(316,17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(316,18).Value = "Región de O'Higgins"
$ws.Cells.Item(316,19).Value = 667
$ws.Cells.Item(316,20).Value = 18

# Row 317 (updated: new quality/price data)
$ws.Cells.Item(317,1).Value = 10
$ws.Cells.Item(317,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(317,3).Value = "La Araucanía"
$ws.Cells.Item(317,4).Value = 44476
$ws.Cells.Item(317,5).Value = 9
$ws.Cells.Item(317,6).Value = "Fruta"
$ws.Cells.Item(317,7).Value = 100102
$ws.Cells.Item(317,8).Value = "Cítricos"
$ws.Cells.Item(317,9).Value = 100102004
$ws.Cells.Item(317,10).Value = "Mandarina"
$ws.Cells.Item(317,11).Value = "Murcott"
$ws.Cells.Item(317,12).Value = "Segunda"
$ws.Cells.Item(317,13).Value = 300
$ws.Cells.Item(317,14).Value = 8000
$ws.Cells.Item(317,15).Value = 8000
$ws.Cells.Item(317,16).Value = 8000
$ws.Cells.Item(317,17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(317,18).Value = "Región de O'Higgins"
$ws.Cells.Item(317,19).Value = 444
$ws.Cells.Item(317,20).Value = 18

# Row 318 (updated: new quality/price data)
$ws.Cells.Item(318,1).Value = 10
$ws.Cells.Item(318,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(318,3).Value = "La Araucanía"
$ws.Cells.Item(318,4).Value = 44476
$ws.Cells.Item(318,5).Value = 9
$ws.Cells.Item(318,6).Value = "Fruta"
$ws.Cells.Item(318,7).Value = 100102
$ws.Cells.Item(318,8).Value = "Cítricos"
$ws.Cells.Item(318,9).Value = 100102004
$ws.Cells.Item(318,10).Value = "Mandarina"
$ws.Cells.Item(318,11).Value = "Murcott"
$ws.Cells.Item(318,12).Value = "Tercera"
$ws.Cells.Item(318,13).Value = 100
$ws.Cells.Item(318,14).Value = 5000
$ws.Cells.Item(318,15).Value = 5000
$ws.Cells.Item(318,16).Value = 5000
$ws.Cells.Item(318,17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(318,18).Value = "Región de O'Higgins"
$ws.Cells.Item(318,19).Value = 278
$ws.Cells.Item(318,20).Value = 18
